$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$r = $ws.Range("D2")
$r.NumberFormat = "@"
$r.Value = '29.434.26'
$r.Style = "Normal"
$ws.Range("E2").Value = '  +0.39%  '
# Row 3
$r = $ws.Range("D3")
$r.NumberFormat = "@"
$r.Value = '1.848.69'
$r.Style = "Normal"
$ws.Range("E3").Value = '  +0.37%  '
# Row 4
$r = $ws.Range("D4")
$r.NumberFormat = "@"
$r.Value = '0.9999'
$r.Style = "Normal"
$ws.Range("E4").Value = '  +0.10%  '
# Row 5
$ws.Range("E5").Value = '  +0.88%  '
# Row 6
$r = $ws.Range("D6")
$r.NumberFormat = "@"
$r.Value = '0.6296'
$r.Style = "Normal"
$ws.Range("E6").Value = '  -0.06%  '
# Row 7
$ws.Range("E7").Value = '  +0.05%  '
# Row 8
$r = $ws.Range("D8")
$r.NumberFormat = "@"
$r.Value = '0.07679'
$r.Style = "Normal"
$ws.Range("E8").Value = '  +2.09%  '
# Row 9
$r = $ws.Range("D9")
$r.NumberFormat = "@"
$r.Value = '0.2924'
$r.Style = "Normal"
# Row 10
$ws.Range("E10").Value = '  +0.97%  '
# Row 11
$r = $ws.Range("D11")
$r.NumberFormat = "@"
$r.Value = '0.07739'
$r.Style = "Normal"
$ws.Range("E11").Value = '  +0.61%  '
# Row 12
$r = $ws.Range("D12")
$r.NumberFormat = "@"
$r.Value = '1.863.64'
$r.Style = "Normal"
$ws.Range("E12").Value = '  +1.41%  '
# Row 13
$ws.Range("E13").Value = '  +1.01%  '
# Row 14
$r = $ws.Range("D14")
$r.NumberFormat = "@"
$r.Value = '0.6796'
$r.Style = "Normal"
$ws.Range("E14").Value = '  +0.25%  '
# Row 15
$r = $ws.Range("D15")
$r.NumberFormat = "@"
$r.Value = '0.00001072'
$r.Style = "Normal"
$ws.Range("E15").Value = '  +2.62%  '
# Row 16
$r = $ws.Range("D16")
$r.NumberFormat = "@"
$r.Value = '83.59'
$r.Style = "Normal"
$ws.Range("E16").Value = '  +0.76%  '
# Row 17
$ws.Range("B17").Value = 'Uniswap'
$ws.Range("C17").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$r = $ws.Range("D17")
$r.NumberFormat = "@"
$r.Value = '6.189'
$r.Style = "Normal"
$ws.Range("E17").Value = '  +0.94%  '
# Row 18
$ws.Range("B18").Value = 'WrappedBTC'
$ws.Range("C18").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$r = $ws.Range("D18")
$r.NumberFormat = "@"
$r.Value = '29.462.48'
$r.Style = "Normal"
$ws.Range("E18").Value = '  +0.50%  '
# Row 19
$ws.Range("B19").Value = 'BitcoinCash'
$ws.Range("C19").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$r = $ws.Range("D19")
$r.NumberFormat = "@"
$r.Value = '228.12'
$r.Style = "Normal"
$ws.Range("E19").Value = '  -0.01%  '
# Row 20
$ws.Range("B20").Value = 'Avalanche'
$ws.Range("C20").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$r = $ws.Range("D20")
$r.NumberFormat = "@"
$r.Value = '12.43'
$r.Style = "Normal"
$ws.Range("E20").Value = '  +0.11%  '
# Row 21
$ws.Range("B21").Value = 'Dai'
$ws.Range("C21").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$r = $ws.Range("D21")
$r.NumberFormat = "@"
$r.Value = '1.000'
$r.Style = "Normal"
$ws.Range("E21").Value = '  +0.06%  '
# Row 22
$ws.Range("B22").Value = 'Chainlink'
$ws.Range("C22").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$r = $ws.Range("D22")
$r.NumberFormat = "@"
$r.Value = '7.428'
$r.Style = "Normal"
$ws.Range("E22").Value = '  +0.21%  '
# Row 23
$ws.Range("B23").Value = 'BinanceUSD'
$ws.Range("C23").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$r = $ws.Range("D23")
$r.NumberFormat = "@"
$r.Value = '1.001'
$r.Style = "Normal"
$ws.Range("E23").Value = '  +0.06%  '
# Row 24
$ws.Range("B24").Value = 'Monero'
$ws.Range("C24").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$r = $ws.Range("D24")
$r.NumberFormat = "@"
$r.Value = '157.96'
$r.Style = "Normal"
$ws.Range("E24").Value = '  +0.87%  '
# Row 25
$ws.Range("B25").Value = 'Stellar'
$ws.Range("C25").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$r = $ws.Range("D25")
$r.NumberFormat = "@"
$r.Value = '0.1380'
$r.Style = "Normal"
$ws.Range("E25").Value = '  -0.63%  '
# Row 26
$ws.Range("B26").Value = 'Cosmos'
$ws.Range("C26").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$r = $ws.Range("D26")
$r.NumberFormat = "@"
$r.Value = '8.406'
$r.Style = "Normal"
$ws.Range("E26").Value = '  +0.77%  '
# Row 27
$ws.Range("B27").Value = 'EthereumClassic'
$ws.Range("C27").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$r = $ws.Range("D27")
$r.NumberFormat = "@"
$r.Value = '17.69'
$r.Style = "Normal"
$ws.Range("E27").Value = '  +0.57%  '
# Row 28
$ws.Range("B28").Value = 'Toncoin'
$ws.Range("C28").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$r = $ws.Range("D28")
$r.NumberFormat = "@"
$r.Value = '1.354'
$r.Style = "Normal"
$ws.Range("E28").Value = '  +6.51%  '
# Row 29
$ws.Range("B29").Value = 'PancakeSwap'
$ws.Range("C29").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$r = $ws.Range("D29")
$r.NumberFormat = "@"
$r.Value = '1.466'
$r.Style = "Normal"
$ws.Range("E29").Value = '  +0.82%  '
# Row 30
$ws.Range("B30").Value = 'Hedera'
$ws.Range("C30").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$r = $ws.Range("D30")
$r.NumberFormat = "@"
$r.Value = '0.05672'
$r.Style = "Normal"
$ws.Range("E30").Value = '  +0.72%  '
# Row 31
$ws.Range("B31").Value = 'Filecoin'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$r = $ws.Range("D31")
$r.NumberFormat = "@"
$r.Value = '4.119'
$r.Style = "Normal"
$ws.Range("E31").Value = '  +0.44%  '
# Row 32
$ws.Range("B32").Value = 'InternetComputer(DFINITY)'
$ws.Range("C32").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$r = $ws.Range("D32")
$r.NumberFormat = "@"
$r.Value = '4.031'
$r.Style = "Normal"
$ws.Range("E32").Value = '  +0.34%  '
# Row 33
$ws.Range("B33").Value = 'LidoDAOToken'
$ws.Range("C33").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$r = $ws.Range("D33")
$r.NumberFormat = "@"
$r.Value = '1.842'
$r.Style = "Normal"
$ws.Range("E33").Value = '  +0.70%  '
# Row 34
$ws.Range("B34").Value = 'ARBITRUM'
$ws.Range("C34").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$r = $ws.Range("D34")
$r.NumberFormat = "@"
$r.Value = '1.161'
$r.Style = "Normal"
$ws.Range("E34").Value = '  +0.61%  '
# Row 35
$ws.Range("B35").Value = 'ImmutableX'
$ws.Range("C35").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$r = $ws.Range("D35")
$r.NumberFormat = "@"
$r.Value = '0.7081'
$r.Style = "Normal"
$ws.Range("E35").Value = '  -0.04%  '
# Row 36
$ws.Range("B36").Value = 'HuobiToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$r = $ws.Range("D36")
$r.NumberFormat = "@"
$r.Value = '2.586'
$r.Style = "Normal"
$ws.Range("E36").Value = '  -0.17%  '
# Row 37
$ws.Range("B37").Value = 'MXToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$r = $ws.Range("D37")
$r.NumberFormat = "@"
$r.Value = '2.778'
$r.Style = "Normal"
$ws.Range("E37").Value = '  +0.72%  '
# Row 38
$ws.Range("B38").Value = 'VeChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$r = $ws.Range("D38")
$r.NumberFormat = "@"
$r.Value = '0.01790'
$r.Style = "Normal"
$ws.Range("E38").Value = '  -1.07%  '
# Row 39
$ws.Range("B39").Value = 'Maker'
$ws.Range("C39").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$r = $ws.Range("D39")
$r.NumberFormat = "@"
$r.Value = '1.220.40'
$r.Style = "Normal"
$ws.Range("E39").Value = '  -1.97%  '
# Row 40
$ws.Range("B40").Value = 'FraxShare'
$ws.Range("C40").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$r = $ws.Range("D40")
$r.NumberFormat = "@"
$r.Value = '6.542'
$r.Style = "Normal"
$ws.Range("E40").Value = '  +4.87%  '
# Row 41
$ws.Range("B41").Value = 'TrustWalletToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$r = $ws.Range("D41")
$r.NumberFormat = "@"
$r.Value = '0.9124'
$r.Style = "Normal"
$ws.Range("E41").Value = '  +1.14%  '
# Row 42
$ws.Range("B42").Value = 'PaxDollar'
$ws.Range("C42").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$r = $ws.Range("D42")
$r.NumberFormat = "@"
$r.Value = '1.001'
$r.Style = "Normal"
$ws.Range("E42").Value = '  +0.16%  '
# Row 43
$ws.Range("B43").Value = 'Quant'
$ws.Range("C43").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$r = $ws.Range("D43")
$r.NumberFormat = "@"
$r.Value = '101.69'
$r.Style = "Normal"
$ws.Range("E43").Value = '  -0.15%  '
# Row 44
$ws.Range("B44").Value = 'Aave'
$ws.Range("C44").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$r = $ws.Range("D44")
$r.NumberFormat = "@"
$r.Value = '66.03'
$r.Style = "Normal"
$ws.Range("E44").Value = '  +0.70%  '
# Row 45
$ws.Range("B45").Value = 'BabyDogeCoin'
$ws.Range("C45").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$r = $ws.Range("D45")
$r.NumberFormat = "@"
$r.Value = '0.00000000121'
$r.Style = "Normal"
$ws.Range("E45").Value = '  +0.76%  '
# Row 46
$ws.Range("B46").Value = 'Aptos'
$ws.Range("C46").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$r = $ws.Range("D46")
$r.NumberFormat = "@"
$r.Value = '7.140'
$r.Style = "Normal"
$ws.Range("E46").Value = '  +0.65%  '
# Row 47
$ws.Range("B47").Value = 'TheSandbox'
$ws.Range("C47").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$r = $ws.Range("D47")
$r.NumberFormat = "@"
$r.Value = '0.4022'
$r.Style = "Normal"
$ws.Range("E47").Value = '  +0.82%  '
# Row 48
$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$r = $ws.Range("D48")
$r.NumberFormat = "@"
$r.Value = '9.049'
$r.Style = "Normal"
$ws.Range("E48").Value = '  +1.84%  '
# Row 49
$ws.Range("B49").Value = 'RenderToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$r = $ws.Range("D49")
$r.NumberFormat = "@"
$r.Value = '1.678'
$r.Style = "Normal"
$ws.Range("E49").Value = '  +0.51%  '
# Row 50
$ws.Range("B50").Value = 'Algorand'
$ws.Range("C50").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$r = $ws.Range("D50")
$r.NumberFormat = "@"
$r.Value = '0.1145'
$r.Style = "Normal"
$ws.Range("E50").Value = '  +2.20%  '
# Row 51
$ws.Range("B51").Value = 'Cronos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$r = $ws.Range("D51")
$r.NumberFormat = "@"
$r.Value = '0.05719'
$r.Style = "Normal"
$ws.Range("E51").Value = '  +0.12%  '
